# Update the generated IG output spreadsheet for the 2025-08 refresh:
#  - new canonical base URL (2rdoc.pt instead of github.com/RicardoLSantos/shorthand)
#  - new generation date/time
# After the text changes, re-run AutoFit on the "Elements" sheet so the
# "best fit" column widths recorded in the worksheet match the regenerated
# content (URLs got shorter, which narrows a few columns).

$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

$oldSd = "https://github.com/RicardoLSantos/shorthand/StructureDefinition/measurement-conditions"
$newSd = "https://2rdoc.pt/ig/ios-lifestyle-medicine/StructureDefinition/measurement-conditions"

$oldVs = "https://github.com/RicardoLSantos/shorthand/ValueSet/measurement-conditions-vs"
$newVs = "https://2rdoc.pt/ig/ios-lifestyle-medicine/ValueSet/measurement-conditions-vs"

$oldDate = "2025-01-24T12:30:11+00:00"
$newDate = "2025-08-20T10:40:04+01:00"

# Metadata sheet: URL (B2) and Date (B8)
$metadata.Range("B2").Value = $newSd
$metadata.Range("B8").Value = $newDate

# Elements sheet: Fixed Value of Extension.url (R5) and Binding Value Set of
# Extension.value[x] (Z6)
$elements.Range("R5").Value = $newSd
$elements.Range("Z6").Value = $newVs

# Recompute the "best fit" column widths on the Elements sheet now that the
# text content has changed.
$elements.Columns.AutoFit() | Out-Null
